# 8.6.1 workbook update:
#  - Fix the Kyrgyz indicator title in A1 (grammar/whitespace correction); the
#    old mistranslated shared string becomes unused and is dropped on save,
#    the new one is appended - this also naturally renumbers/repacks the
#    other shared strings (Kyrgyz Republic / Kyrgyz Republic(ky) / Kyrgyz
#    Republic(en) in row 5) without touching them directly.
#  - A1's vertical alignment moves from "top" to "center" (still left/wrap).
#  - Row 1 height shrinks from 54 to 48.
#  - A new 2023 data column (T) is added, mirroring the formatting already
#    used for 2022 (column S) in the header/year row and the three data rows.
#  - Reset the lingering "S4:S7" selection left over in the sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- A1: corrected Kyrgyz title text + vertical-center alignment ---------
$ws.Range("A1").Value = "8.6.1 Иштебеген, окубаган жана кесиптик көндүмдөрдү үйрөнбөгөн (15 жаштан 24 жашка чейинки ) жаштардын үлүшү  "
$ws.Range("A1").VerticalAlignment = -4108
$ws.Rows.Item(1).RowHeight = 48

# --- New column T: 2023 ----------------------------------------------------
# Copy number/font/border formatting from the existing 2022 column (S) so the
# new column matches its neighbours, then overwrite the values.
$ws.Range("S4:S7").Copy()
$ws.Range("T4:T7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("T4").Value = 2023
$ws.Range("T5").Value = 18.6
$ws.Range("T6").Value = 11.5
$ws.Range("T7").Value = 25.9

# --- Clear the stale S4:S7 selection left in the saved view ---------------
$ws.Range("A1").Select()
